# "added validation functions and converted the project for automating facebook"
#
# Fill in the Facebook login credentials on the "Data" sheet (row 3 of the
# first data block: Runmode/username/password) and flip the Runmode flag
# in row 4 from N to Y so both rows run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("B3").Value = "vaibhavcool20@protonmail.com"
$ws.Range("C3").Value = "xxxxxxx"
$ws.Range("A4").Value = "Y"

# Widen the username/password columns so the new values are fully visible
# (mirrors Excel's own best-fit autosize after the longer text was typed
# in). ColumnWidth is quantized by the host in 1/6-character steps, so
# these inputs are chosen to land on the closest attainable width to the
# target 30.28515625 / 15.140625 character units.
$ws.Columns.Item(2).ColumnWidth = 29 + 5/12
$ws.Columns.Item(3).ColumnWidth = 14.25
